$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new J3 = "Indirect" addressing mode label (new shared string #37)
$ws.Range("J3").Value = "Indirect"

# Row 22: instruction description for the new Indirect instruction (new shared string #38)
$ws.Range("J22").Value = "Move to address R14 (52) data in R8 (72)"

# Row 2: J2 becomes the text "4DE8" (was numeric 4090) (new shared string #39)
$ws.Range("J2").Value = "4DE8"

# New "Indirect" column mirrors column G values for rows 6-20
$ws.Range("J6").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("J9").Value = 20
$ws.Range("J10").Value = 35
$ws.Range("J11").Value = 41
$ws.Range("J12").Value = 55
$ws.Range("J13").Value = 72
$ws.Range("J14").Value = 50
$ws.Range("J15").Value = 51
$ws.Range("J16").Value = 422
$ws.Range("J17").Value = 874
$ws.Range("J18").Value = 52
$ws.Range("J19").Value = 53
$ws.Range("J20").Value = -4

# Row 22 grows taller to fit the new wrapped "Indirect" description
$ws.Rows.Item(22).RowHeight = 75

# Update selection to match the author's final cursor position
$ws.Range("J2").Select()
